$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("B3").Value = 1.02
$ws.Range("B4").Value = 1.02
$ws.Range("B5").Value = 1.02
$ws.Range("B6").Value = 1.02
$ws.Range("B7").Value = 1.02
$ws.Range("B8").Value = 1.02
$ws.Range("B9").Value = 1.02
$ws.Range("B10").Value = 1.02
$ws.Range("B11").Value = 1.02
$ws.Range("B12").Value = 1.02
$ws.Range("B13").Value = 1.02
$ws.Range("B14").Value = 1.02
$ws.Range("B15").Value = 1.02
$ws.Range("B16").Value = 1.02
$ws.Range("B17").Value = 1.02
$ws.Range("B18").Value = 1.02
$ws.Range("B19").Value = 1.02
$ws.Range("B20").Value = 1.02
$ws.Range("B21").Value = 1.02
$ws.Range("B22").Value = 1.02
$ws.Range("B23").Value = 1.02
$ws.Range("B24").Value = 1.02
$ws.Range("B25").Value = 1.02

$ws.Range("C2").Value = 1.045651583328221
$ws.Range("C3").Value = 1.046595137350151
$ws.Range("C4").Value = 1.047206259149114
$ws.Range("C5").Value = 1.047463312196275
$ws.Range("C6").Value = 1.047506480557994
$ws.Range("C7").Value = 1.047209693365989
$ws.Range("C8").Value = 1.045970341147802
$ws.Range("C9").Value = 1.043790935043896
$ws.Range("C10").Value = 1.042341089540228
$ws.Range("C11").Value = 1.04171403877103
$ws.Range("C12").Value = 1.041481236826643
$ws.Range("C13").Value = 1.04153116851839
$ws.Range("C14").Value = 1.041694792985275
$ws.Range("C15").Value = 1.041795622379316
$ws.Range("C16").Value = 1.042382720565268
$ws.Range("C17").Value = 1.042751191197797
$ws.Range("C18").Value = 1.042966185387824
$ws.Range("C19").Value = 1.043039504943538
$ws.Range("C20").Value = 1.04271165038163
$ws.Range("C21").Value = 1.041646606539442
$ws.Range("C22").Value = 1.040977623433141
$ws.Range("C23").Value = 1.041332201820153
$ws.Range("C24").Value = 1.042729516953601
$ws.Range("C25").Value = 1.044353824118709

$ws.Range("D2").Value = 1.053044234695989
$ws.Range("D3").Value = 1.053878671815083
$ws.Range("D4").Value = 1.054419437638948
$ws.Range("D5").Value = 1.054646972546816
$ws.Range("D6").Value = 1.054685188163294
$ws.Range("D7").Value = 1.05442247719958
$ws.Range("D8").Value = 1.053326063996193
$ws.Range("D9").Value = 1.051400463497025
$ws.Range("D10").Value = 1.050121138495623
$ws.Range("D11").Value = 1.049568241690959
$ws.Range("D12").Value = 1.049363031907343
$ws.Range("D13").Value = 1.049407042805294
$ws.Range("D14").Value = 1.049551275698883
$ws.Range("D15").Value = 1.049640163686562
$ws.Range("D16").Value = 1.050157854887496
$ws.Range("D17").Value = 1.05048287361643
$ws.Range("D18").Value = 1.050672553651364
$ws.Range("D19").Value = 1.050737246913848
$ws.Range("D20").Value = 1.050447991609921
$ws.Range("D21").Value = 1.049508798257872
$ws.Range("D22").Value = 1.04891922069166
$ws.Range("D23").Value = 1.04923167818972
$ws.Range("D24").Value = 1.050463752971854
$ws.Range("D25").Value = 1.051897507094739

$ws.Range("E2").Value = 1.053101525993051
$ws.Range("E3").Value = 1.053932154318194
$ws.Range("E4").Value = 1.054470499162553
$ws.Range("E5").Value = 1.054697026713693
$ws.Range("E6").Value = 1.054735073801566
$ws.Range("E7").Value = 1.054473525221888
$ws.Range("E8").Value = 1.053382059069266
$ws.Range("E9").Value = 1.051465509845913
$ws.Range("E10").Value = 1.050192443962306
$ws.Range("E11").Value = 1.04964231087225
$ws.Range("E12").Value = 1.049438135698852
$ws.Range("E13").Value = 1.049481924305261
$ws.Range("E14").Value = 1.049625430237237
$ws.Range("E15").Value = 1.049713871386004
$ws.Range("E16").Value = 1.050228978063021
$ws.Range("E17").Value = 1.050552389911998
$ws.Range("E18").Value = 1.050741137840004
$ws.Range("E19").Value = 1.050805514152871
$ws.Range("E20").Value = 1.050517679774818
$ws.Range("E21").Value = 1.049583166646185
$ws.Range("E22").Value = 1.048996578264312
$ws.Range("E23").Value = 1.049307446724415
$ws.Range("E24").Value = 1.050533363460462
$ws.Range("E25").Value = 1.051960173855275

$ws.Range("F2").Value = 1.063263684837116
$ws.Range("F3").Value = 1.06415665900342
$ws.Range("F4").Value = 1.064735399182342
$ws.Range("F5").Value = 1.064978921651741
$ws.Range("F6").Value = 1.065019823005148
$ws.Range("F7").Value = 1.064738652278709
$ws.Range("F8").Value = 1.063565276884279
$ws.Range("F9").Value = 1.061504805070521
$ws.Range("F10").Value = 1.060136071667699
$ws.Range("F11").Value = 1.059544580908772
$ws.Range("F12").Value = 1.05932505378376
$ws.Range("F13").Value = 1.059372134976581
$ws.Range("F14").Value = 1.059526431069678
$ws.Range("F15").Value = 1.059621521736772
$ws.Range("F16").Value = 1.060175351957571
$ws.Range("F17").Value = 1.060523072043648
$ws.Range("F18").Value = 1.06072600501524
$ws.Range("F19").Value = 1.060795219155262
$ws.Range("F20").Value = 1.060485753194023
$ws.Range("F21").Value = 1.059480989765027
$ws.Range("F22").Value = 1.058850290860166
$ws.Range("F23").Value = 1.059184537564779
$ws.Range("F24").Value = 1.060502615624575
$ws.Range("F25").Value = 1.062036627300648

$ws.Range("I2").Value = 1.038196263514373
$ws.Range("I3").Value = 1.038320663966592
$ws.Range("I4").Value = 1.038399867441055
$ws.Range("I5").Value = 1.038432854714927
$ws.Range("I6").Value = 1.038438375235941
$ws.Range("I7").Value = 1.038400309436903
$ws.Range("I8").Value = 1.038238572387942
$ws.Range("I9").Value = 1.037943701752968
$ws.Range("I10").Value = 1.037740519852296
$ws.Range("I11").Value = 1.037650981177011
$ws.Range("I12").Value = 1.037617488766639
$ws.Range("I13").Value = 1.037624683571658
$ws.Range("I14").Value = 1.037648217449372
$ws.Range("I15").Value = 1.037662686489147
$ws.Range("I16").Value = 1.037746429427683
$ws.Range("I17").Value = 1.037798542005704
$ws.Range("I18").Value = 1.037828787877916
$ws.Range("I19").Value = 1.037839075392459
$ws.Range("I20").Value = 1.037792966383561
$ws.Range("I21").Value = 1.037641293757453
$ws.Range("I22").Value = 1.037544578968558
$ws.Range("I23").Value = 1.037595977296536
$ws.Range("I24").Value = 1.037795486232324
$ws.Range("I25").Value = 1.038021098934686

$ws.Range("J2").Value = 1.050710177185332
$ws.Range("J3").Value = 1.051301509324517
$ws.Range("J4").Value = 1.051684082784484
$ws.Range("J5").Value = 1.051844901631194
$ws.Range("J6").Value = 1.051871902907884
$ws.Range("J7").Value = 1.051686231714316
$ws.Range("J8").Value = 1.050910031959031
$ws.Range("J9").Value = 1.049541875887237
$ws.Range("J10").Value = 1.048629576091556
$ws.Range("J11").Value = 1.04823450806331
$ws.Range("J12").Value = 1.048087757897876
$ws.Range("J13").Value = 1.048119236483389
$ws.Range("J14").Value = 1.048222377730234
$ws.Range("J15").Value = 1.04828592591926
$ws.Range("J16").Value = 1.048655794789171
$ws.Range("J17").Value = 1.048887794865171
$ws.Range("J18").Value = 1.049023113016175
$ws.Range("J19").Value = 1.049069252369183
$ws.Range("J20").Value = 1.048862903813553
$ws.Range("J21").Value = 1.048192005313654
$ws.Range("J22").Value = 1.047770160425961
$ws.Range("J23").Value = 1.047993790315569
$ws.Range("J24").Value = 1.048874151019495
$ws.Range("J25").Value = 1.049895616542289

$ws.Range("K2").Value = 1.055791291273662
$ws.Range("K3").Value = 1.056438639673142
$ws.Range("K4").Value = 1.056857696862197
$ws.Range("K5").Value = 1.057033910214919
$ws.Range("K6").Value = 1.057063499621544
$ws.Range("K7").Value = 1.056860051271327
$ws.Range("K8").Value = 1.056010027774774
$ws.Range("K9").Value = 1.054513608063746
$ws.Range("K10").Value = 1.05351703115504
$ws.Range("K11").Value = 1.053085763722236
$ws.Range("K12").Value = 1.052925611636508
$ws.Range("K13").Value = 1.052959962979774
$ws.Range("K14").Value = 1.053072524680964
$ws.Range("K15").Value = 1.053141882976433
$ws.Range("K16").Value = 1.053545658473395
$ws.Range("K17").Value = 1.053799005840206
$ws.Range("K18").Value = 1.05394680378138
$ws.Range("K19").Value = 1.053997203213206
$ws.Range("K20").Value = 1.053771821500744
$ws.Range("K21").Value = 1.053039376955738
$ws.Range("K22").Value = 1.052579091120735
$ws.Range("K23").Value = 1.05282307502961
$ws.Range("K24").Value = 1.053784104857395
$ws.Range("K25").Value = 1.054900291247993

$ws.Range("L2").Value = 1.055848424269224
$ws.Range("L3").Value = 1.05649198507054
$ws.Range("L4").Value = 1.056908634070029
$ws.Range("L5").Value = 1.0570838452296
$ws.Range("L6").Value = 1.057113266962336
$ws.Range("L7").Value = 1.056910975047686
$ws.Range("L8").Value = 1.056065871906025
$ws.Range("L9").Value = 1.054578448982281
$ws.Range("L10").Value = 1.053588089363801
$ws.Range("L11").Value = 1.053159566121891
$ws.Range("L12").Value = 1.053000441173297
$ws.Range("L13").Value = 1.053034571838171
$ws.Range("L14").Value = 1.053146411824516
$ws.Range("L15").Value = 1.053215326484214
$ws.Range("L16").Value = 1.053616535656206
$ws.Range("L17").Value = 1.053868287174386
$ws.Range("L18").Value = 1.054015159307368
$ws.Range("L19").Value = 1.054065243914549
$ws.Range("L20").Value = 1.053841273534791
$ws.Range("L21").Value = 1.053113476410382
$ws.Range("L22").Value = 1.052656157851627
$ws.Range("L23").Value = 1.052898564461095
$ws.Range("L24").Value = 1.053853479744
$ws.Range("L25").Value = 1.054962767598397

$ws.Range("M2").Value = 1.065982796224523
$ws.Range("M3").Value = 1.066690551378756
$ws.Range("M4").Value = 1.067148802977416
$ws.Range("M5").Value = 1.067341519203886
$ws.Range("M6").Value = 1.067373881010051
$ws.Range("M7").Value = 1.067151377798534
$ws.Range("M8").Value = 1.066221925202636
$ws.Range("M9").Value = 1.064586367599357
$ws.Range("M10").Value = 1.063497594822893
$ws.Range("M11").Value = 1.063026540540886
$ws.Range("M12").Value = 1.062851630187843
$ws.Range("M13").Value = 1.062889146307153
$ws.Range("M14").Value = 1.063012081160996
$ws.Range("M15").Value = 1.063087833397077
$ws.Range("M16").Value = 1.063528865509076
$ws.Range("M17").Value = 1.06380561908078
$ws.Range("M18").Value = 1.063967082441108
$ws.Range("M19").Value = 1.064022143624923
$ws.Range("M20").Value = 1.063775922134012
$ws.Range("M21").Value = 1.062975878238392
$ws.Range("M22").Value = 1.06247320775012
$ws.Range("M23").Value = 1.062739649364246
$ws.Range("M24").Value = 1.063789340789672
$ws.Range("M25").Value = 1.065008922602883

$ws.Range("N2").Value = 1.052202306155078
$ws.Range("N3").Value = 1.052794478053727
$ws.Range("N4").Value = 1.053177594811889
$ws.Range("N5").Value = 1.053338642039807
$ws.Range("N6").Value = 1.053365681661407
$ws.Range("N7").Value = 1.053179746793447
$ws.Range("N8").Value = 1.052402444745479
$ws.Range("N9").Value = 1.051032345735133
$ws.Range("N10").Value = 1.050118750369099
$ws.Range("N11").Value = 1.049723121298939
$ws.Range("N12").Value = 1.049576162731438
$ws.Range("N13").Value = 1.049607686020152
$ws.Range("N14").Value = 1.049710973739398
$ws.Range("N15").Value = 1.049774612174141
$ws.Range("N16").Value = 1.050145006300273
$ws.Range("N17").Value = 1.050377335842991
$ws.Range("N18").Value = 1.050512846161291
$ws.Range("N19").Value = 1.050559051037472
$ws.Range("N20").Value = 1.050352409443225
$ws.Range("N21").Value = 1.049680558190503
$ws.Range("N22").Value = 1.049258114234683
$ws.Range("N23").Value = 1.049482061704385
$ws.Range("N24").Value = 1.05036367262149
$ws.Range("N25").Value = 1.05138658874249
